$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 - copy style from H1 (bold, bordered, centered)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Data values for columns I (I0) and J (IF), rows 2-43
$data = @(
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(7, 8),
    @(6, 7),
    @(8, 8),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(7, 7),
    @(9, 9),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(7, 7),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(7, 8),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(8, 9),
    @(9, 9),
    @(8, 8),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(11, 11),
    @(8, 8),
    @(9, 9),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(4, 4)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
